$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trial-type reordering: the stimulus (col A) and the per-key i/e mapping
# (cols D:K) for rows 17-31 get cyclically rotated down by 5 rows within
# that 15-row block (17-31), while columns B (block #) and C (colour) stay
# put since they describe the physical row/trial slot.
#
# i.e. new(row) = old(row - 5), wrapping row-5 < 17 to row+10.

$firstRow = 17
$lastRow = 31
$n = $lastRow - $firstRow + 1  # 15

# Snapshot the values that need to move, before we overwrite anything.
$oldA = $ws.Range("A$firstRow`:A$lastRow").Value2
$oldDK = $ws.Range("D$firstRow`:K$lastRow").Value2

$newA = New-Object 'object[,]' $n,1
$newDK = New-Object 'object[,]' $n,8

for ($i = 1; $i -le $n; $i++) {
    $srcIdx = $i - 5
    if ($srcIdx -lt 1) { $srcIdx = $srcIdx + $n }

    $newA[$i-1,0] = $oldA[$srcIdx,1]
    for ($j = 1; $j -le 8; $j++) {
        $newDK[$i-1,$j-1] = $oldDK[$srcIdx,$j]
    }
}

$ws.Range("A$firstRow`:A$lastRow").Value2 = $newA
$ws.Range("D$firstRow`:K$lastRow").Value2 = $newDK

$ws.Range("B32").Select() | Out-Null
